$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2 through 42) holds a date serial value (45710) that
# needs to be incremented by one day (45711), keeping existing formatting.
for ($row = 2; $row -le 42; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45710) {
        $cell.Value2 = 45711
    }
}
